$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.091.57"
$ws.Range("E2").Value = "  -1.41%  "

$ws.Range("D3").Value = "1.833.97"
$ws.Range("E3").Value = "  -1.18%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'240.15"
$ws.Range("E5").Value = "  -2.11%  "

$ws.Range("D6").Value = "'0.6607"
$ws.Range("E6").Value = "  -4.80%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "'0.2940"
$ws.Range("E8").Value = "  -4.05%  "

$ws.Range("D9").Value = "'0.07343"
$ws.Range("E9").Value = "  -4.19%  "

$ws.Range("D10").Value = "'22.67"
$ws.Range("E10").Value = "  -3.61%  "

$ws.Range("D11").Value = "'0.07663"
$ws.Range("E11").Value = "  -1.26%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'5.007"
$ws.Range("E12").Value = "  -2.56%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.6731"
$ws.Range("E13").Value = "  -2.91%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.733.53"
$ws.Range("E14").Value = "  -6.61%  "

$ws.Range("D15").Value = "'85.88"
$ws.Range("E15").Value = "  -5.38%  "

$ws.Range("D16").Value = "'6.168"
$ws.Range("E16").Value = "  -2.04%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.000008198"
$ws.Range("E17").Value = "  -1.00%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "28.694.15"
$ws.Range("E18").Value = "  -2.82%  "

$ws.Range("D19").Value = "'227.79"
$ws.Range("E19").Value = "  -3.64%  "

$ws.Range("D20").Value = "'12.47"
$ws.Range("E20").Value = "  -1.88%  "

$ws.Range("D21").Value = "'0.9996"
$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("D22").Value = "'7.263"
$ws.Range("E22").Value = "  -4.66%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").Value = "'160.60"
$ws.Range("E24").Value = "  +0.66%  "

$ws.Range("D25").Value = "'8.673"
$ws.Range("E25").Value = "  -2.61%  "

$ws.Range("D26").Value = "'0.1402"
$ws.Range("E26").Value = "  -5.44%  "

$ws.Range("D27").Value = "'18.00"
$ws.Range("E27").Value = "  -1.20%  "

$ws.Range("D28").Value = "'1.502"
$ws.Range("E28").Value = "  -2.11%  "

$ws.Range("D29").Value = "'4.205"
$ws.Range("E29").Value = "  -0.96%  "

$ws.Range("D30").Value = "'4.070"
$ws.Range("E30").Value = "  -1.67%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.05340"
$ws.Range("E31").Value = "  +2.19%  "

$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'1.178"
$ws.Range("E32").Value = "  -1.94%  "

$ws.Range("D33").Value = "'1.849"
$ws.Range("E33").Value = "  -1.28%  "

$ws.Range("D34").Value = "'0.7460"
$ws.Range("E34").Value = "  -4.03%  "

$ws.Range("D35").Value = "'1.128"
$ws.Range("E35").Value = "  -1.54%  "

$ws.Range("D36").Value = "'2.683"
$ws.Range("E36").Value = "  +0.20%  "

$ws.Range("D37").Value = "1.316.66"
$ws.Range("E37").Value = "  +0.21%  "

$ws.Range("D38").Value = "'0.01802"
$ws.Range("E38").Value = "  -3.52%  "

$ws.Range("D39").Value = "'2.712"
$ws.Range("E39").Value = "  -0.47%  "

$ws.Range("D40").Value = "'0.9148"
$ws.Range("E40").Value = "  -3.08%  "

$ws.Range("D41").Value = "'5.964"
$ws.Range("E41").Value = "  +2.81%  "

$ws.Range("D42").Value = "'0.9985"
$ws.Range("E42").Value = "  -0.22%  "

$ws.Range("D43").Value = "'103.20"
$ws.Range("E43").Value = "  -2.69%  "

$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").Value = "'0.00000000123"
$ws.Range("E44").Value = "  +1.20%  "

$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.5162"
$ws.Range("E45").Value = "  -1.28%  "

$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.932.86"
$ws.Range("E46").Value = "  -3.54%  "

$ws.Range("D47").Value = "'63.69"
$ws.Range("E47").Value = "  +1.44%  "

$ws.Range("D48").Value = "'1.753"
$ws.Range("E48").Value = "  -1.49%  "

$ws.Range("B49").Value = "XinFinNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D49").Value = "'0.07496"
$ws.Range("E49").Value = "  +11.93%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.198"
$ws.Range("E50").Value = "  -5.70%  "

$ws.Range("D51").Value = "'0.05919"
$ws.Range("E51").Value = "  -0.53%  "
